# Updates the cryptocurrency price/volume snapshot on the active sheet to
# match a fresh pull of the feed. Refreshes the Price (D) and Volume(1h)
# (E) columns for most rows; for rows 13/14, 38/39 and 50/51 the
# underlying coins swapped rank order, so Coin (B), Link (C), Price (D)
# and Volume(1h) (E) are all rewritten for that pair.
#
# Price values are plain text in this sheet (not numbers), so any new
# price that Excel would otherwise auto-parse as a number is entered
# with a leading apostrophe to force text, matching the original
# formatting (t="inlineStr").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.438.24"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.794.58"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'226.93"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'32.38"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "'0.0692"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").Value = "'0.0951"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "2.053.15"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.802.20"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.04"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "34.382.21"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "'68.28"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "0.0₃0800"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("D20").Value = "'246.20"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "'11.05"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'4.15"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").Value = "'162.64"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  +2.15%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  +7.85%  "
$ws.Range("E33").Value = "  +3.14%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").Value = "1.440.87"
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("D36").Value = "'2.62"
$ws.Range("E36").Value = "  +7.75%  "
$ws.Range("D37").Value = "'0.671"
$ws.Range("E37").Value = "  +3.94%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0191"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.05"
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("D40").Value = "'84.07"
$ws.Range("E40").Value = "  +4.62%  "
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("D42").Value = "'0.932"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").Value = "'2.74"
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("E45").Value = "  +3.92%  "
$ws.Range("D46").Value = "'6.10"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "1.950.20"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "'105.76"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0130"
$ws.Range("E51").Value = "  -4.81%  "
